$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.425.22"
$ws.Range("E2").Value = "  +4.60%  "
$ws.Range("D3").Value = "3.628.30"
$ws.Range("E3").Value = "  +7.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.45"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.33"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "3.619.17"
$ws.Range("E7").Value = "  +7.13%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.50"
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000290"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "710.07"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "4.214.74"
$ws.Range("E15").Value = "  +7.15%  "
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "72.554.74"
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("D18").Value = "3.603.08"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.934"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.87"
$ws.Range("E23").Value = "  +8.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.84"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "105.49"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("E28").Value = "  +4.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.53"
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.45"
$ws.Range("E31").Value = "  +7.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.19"
$ws.Range("E32").Value = "  +15.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "595.10"
$ws.Range("E33").Value = "  +7.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.36"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.84"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("D39").Value = "3.645.05"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  +8.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.91"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  +7.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0450"
$ws.Range("E44").Value = "  +6.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.351"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("E48").Value = "  +5.24%  "
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.55"
$ws.Range("E51").Value = "  -0.02%  "
